# Refresh cached market-price / Leve-profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
# across all 8 sheets, per the scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2:
$ws.Cells.Item(2, 8).Value = 3263.5715
$ws.Cells.Item(2, 9).Value = 2468.6
$ws.Cells.Item(2, 11).Value = 2468.6
$ws.Cells.Item(2, 13).Value = -2355.6
# Row 4: Root Rush
$ws.Cells.Item(4, 8).Value = 320.64517
$ws.Cells.Item(4, 9).Value = 212.14285
$ws.Cells.Item(4, 10).Value = 1333.3334
$ws.Cells.Item(4, 11).Value = 212.14285
$ws.Cells.Item(4, 12).Value = 1333.3334
$ws.Cells.Item(4, 13).Value = -98.14285000000001
$ws.Cells.Item(4, 14).Value = -1561.3334
# Row 80: Cleansing the Wicked Humours
$ws.Cells.Item(80, 8).Value = 1147.7368
$ws.Cells.Item(80, 9).Value = 275.7143
$ws.Cells.Item(80, 10).Value = 1656.4166
$ws.Cells.Item(80, 11).Value = 827.1428999999999
$ws.Cells.Item(80, 12).Value = 4969.2498
$ws.Cells.Item(80, 13).Value = 170.8571000000001
$ws.Cells.Item(80, 14).Value = -6965.2498
# Row 83: Washing Away the Sins (L)
$ws.Cells.Item(83, 8).Value = 1147.7368
$ws.Cells.Item(83, 9).Value = 275.7143
$ws.Cells.Item(83, 10).Value = 1656.4166
$ws.Cells.Item(83, 11).Value = 2481.4287
$ws.Cells.Item(83, 12).Value = 14907.7494
$ws.Cells.Item(83, 13).Value = 2510.5713
$ws.Cells.Item(83, 14).Value = -24891.7494
# Row 98: The Dotted Line
$ws.Cells.Item(98, 8).Value = 784.9474
$ws.Cells.Item(98, 10).Value = 125
$ws.Cells.Item(98, 12).Value = 125
$ws.Cells.Item(98, 14).Value = -3121
# Row 122: Wishful Inking
$ws.Cells.Item(122, 8).Value = 784.9474
$ws.Cells.Item(122, 10).Value = 125
$ws.Cells.Item(122, 12).Value = 375
$ws.Cells.Item(122, 14).Value = -5275
# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 4409.313
$ws.Cells.Item(138, 9).Value = 3264.4583
$ws.Cells.Item(138, 10).Value = 4775.6665
$ws.Cells.Item(138, 11).Value = 9793.374899999999
$ws.Cells.Item(138, 12).Value = 14326.9995
$ws.Cells.Item(138, 13).Value = -4653.374899999999
$ws.Cells.Item(138, 14).Value = -24606.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 23: A Well-rounded Crew
$ws.Cells.Item(23, 8).Value = 22000
$ws.Cells.Item(23, 10).Value = 22000
$ws.Cells.Item(23, 12).Value = 22000
$ws.Cells.Item(23, 14).Value = -22518
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 17631.164
$ws.Cells.Item(32, 9).Value = 17320.617
$ws.Cells.Item(32, 10).Value = 19856.75
$ws.Cells.Item(32, 11).Value = 17320.617
$ws.Cells.Item(32, 12).Value = 19856.75
$ws.Cells.Item(32, 13).Value = -17033.617
$ws.Cells.Item(32, 14).Value = -20430.75
# Row 41: Skillet Scandal
$ws.Cells.Item(41, 8).Value = 440.25
$ws.Cells.Item(41, 9).Value = 440.25
$ws.Cells.Item(41, 11).Value = 440.25
$ws.Cells.Item(41, 13).Value = -26.25
# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 2607.875
$ws.Cells.Item(45, 9).Value = 1063.3334
$ws.Cells.Item(45, 10).Value = 3534.6
$ws.Cells.Item(45, 11).Value = 1063.3334
$ws.Cells.Item(45, 12).Value = 3534.6
$ws.Cells.Item(45, 13).Value = -686.3334
$ws.Cells.Item(45, 14).Value = -4288.6
# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 418251.78
$ws.Cells.Item(74, 9).Value = 477621.16
$ws.Cells.Item(74, 10).Value = 2666.3333
$ws.Cells.Item(74, 11).Value = 477621.16
$ws.Cells.Item(74, 12).Value = 2666.3333
$ws.Cells.Item(74, 13).Value = -476747.16
$ws.Cells.Item(74, 14).Value = -4414.3333
# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 418251.78
$ws.Cells.Item(77, 9).Value = 477621.16
$ws.Cells.Item(77, 10).Value = 2666.3333
$ws.Cells.Item(77, 11).Value = 2388105.8
$ws.Cells.Item(77, 12).Value = 13331.6665
$ws.Cells.Item(77, 13).Value = -2383737.8
$ws.Cells.Item(77, 14).Value = -22067.6665
# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 4946.2383
$ws.Cells.Item(122, 9).Value = 4634.579
$ws.Cells.Item(122, 11).Value = 13903.737
$ws.Cells.Item(122, 13).Value = -11453.737
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 6899.5093
$ws.Cells.Item(132, 9).Value = 3149.3953
$ws.Cells.Item(132, 10).Value = 23025
$ws.Cells.Item(132, 11).Value = 9448.1859
$ws.Cells.Item(132, 12).Value = 69075
$ws.Cells.Item(132, 13).Value = -6918.1859
$ws.Cells.Item(132, 14).Value = -74135

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector
$ws.Cells.Item(82, 8).Value = 21731.945
$ws.Cells.Item(82, 9).Value = 21731.945
$ws.Cells.Item(82, 11).Value = 21731.945
$ws.Cells.Item(82, 13).Value = -21348.945
# Row 85: The Clamor for Hammers (L)
$ws.Cells.Item(85, 8).Value = 21731.945
$ws.Cells.Item(85, 9).Value = 21731.945
$ws.Cells.Item(85, 11).Value = 21731.945
$ws.Cells.Item(85, 13).Value = -20405.945
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 92783.73
$ws.Cells.Item(86, 9).Value = 1902.4706
$ws.Cells.Item(86, 11).Value = 1902.4706
$ws.Cells.Item(86, 13).Value = -779.4706000000001
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 92783.73
$ws.Cells.Item(89, 9).Value = 1902.4706
$ws.Cells.Item(89, 11).Value = 9512.353000000001
$ws.Cells.Item(89, 13).Value = -3896.353000000001
# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 1249.8077
$ws.Cells.Item(94, 10).Value = 1589.1818
$ws.Cells.Item(94, 12).Value = 1589.1818
$ws.Cells.Item(94, 14).Value = -2491.1818
# Row 97: File under Dull
$ws.Cells.Item(97, 8).Value = 6661
$ws.Cells.Item(97, 9).Value = 6828.75
$ws.Cells.Item(97, 10).Value = 5990
$ws.Cells.Item(97, 11).Value = 6828.75
$ws.Cells.Item(97, 12).Value = 5990
$ws.Cells.Item(97, 13).Value = -5837.75
$ws.Cells.Item(97, 14).Value = -7972

$ws = $wb.Worksheets.Item("CRP")
# Row 94: Beech, Please
$ws.Cells.Item(94, 8).Value = 76434.86
$ws.Cells.Item(94, 9).Value = 131076
$ws.Cells.Item(94, 11).Value = 131076
$ws.Cells.Item(94, 13).Value = -130625
# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 876.55554
$ws.Cells.Item(122, 9).Value = 948.875
$ws.Cells.Item(122, 10).Value = 298
$ws.Cells.Item(122, 11).Value = 2846.625
$ws.Cells.Item(122, 12).Value = 894
$ws.Cells.Item(122, 13).Value = -396.625
$ws.Cells.Item(122, 14).Value = -5794
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 19119
$ws.Cells.Item(132, 9).Value = 2680.6924
$ws.Cells.Item(132, 11).Value = 8042.0772
$ws.Cells.Item(132, 13).Value = -5512.0772
# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 3412.739
$ws.Cells.Item(134, 9).Value = 3009.9092
$ws.Cells.Item(134, 10).Value = 12275
$ws.Cells.Item(134, 11).Value = 9029.7276
$ws.Cells.Item(134, 12).Value = 36825
$ws.Cells.Item(134, 13).Value = -6494.7276
$ws.Cells.Item(134, 14).Value = -41895

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Cells.Item(14, 8).Value = 834.25
$ws.Cells.Item(14, 9).Value = 834.25
$ws.Cells.Item(14, 11).Value = 2502.75
$ws.Cells.Item(14, 13).Value = -2329.75
# Row 68: Such a Butter Face
$ws.Cells.Item(68, 8).Value = 3527.7778
$ws.Cells.Item(68, 9).Value = 975
$ws.Cells.Item(68, 10).Value = 4257.143
$ws.Cells.Item(68, 11).Value = 2925
$ws.Cells.Item(68, 12).Value = 12771.429
$ws.Cells.Item(68, 13).Value = -2114
$ws.Cells.Item(68, 14).Value = -14393.429
# Row 71: No Margarine of Error (L)
$ws.Cells.Item(71, 8).Value = 3527.7778
$ws.Cells.Item(71, 9).Value = 975
$ws.Cells.Item(71, 10).Value = 4257.143
$ws.Cells.Item(71, 11).Value = 8775
$ws.Cells.Item(71, 12).Value = 38314.287
$ws.Cells.Item(71, 13).Value = -4719
$ws.Cells.Item(71, 14).Value = -46426.287
# Row 103: West Meats East
$ws.Cells.Item(103, 8).Value = 328.44446
$ws.Cells.Item(103, 10).Value = 374.42856
$ws.Cells.Item(103, 12).Value = 1123.28568
$ws.Cells.Item(103, 14).Value = -2881.28568
# Row 122: Salt of the North
$ws.Cells.Item(122, 10).Value = 1371.75
$ws.Cells.Item(122, 12).Value = 12345.75
$ws.Cells.Item(122, 14).Value = -17245.75

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 3714
$ws.Cells.Item(80, 9).Value = 2250
$ws.Cells.Item(80, 11).Value = 2250
$ws.Cells.Item(80, 13).Value = -1252
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 3714
$ws.Cells.Item(83, 9).Value = 2250
$ws.Cells.Item(83, 11).Value = 11250
$ws.Cells.Item(83, 13).Value = -6258
# Row 99: Needle in a Hingan Stack
$ws.Cells.Item(99, 8).Value = 14280.125
$ws.Cells.Item(99, 9).Value = 14280.125
$ws.Cells.Item(99, 11).Value = 14280.125
$ws.Cells.Item(99, 13).Value = -12034.125
# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 2769.6667
$ws.Cells.Item(102, 9).Value = 2769.6667
$ws.Cells.Item(102, 11).Value = 2769.6667
$ws.Cells.Item(102, 13).Value = -1147.6667
# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 1796.4445
$ws.Cells.Item(122, 10).Value = 1437.5
$ws.Cells.Item(122, 12).Value = 4312.5
$ws.Cells.Item(122, 14).Value = -9212.5
# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 4858.147
$ws.Cells.Item(132, 9).Value = 3705.2693
$ws.Cells.Item(132, 11).Value = 11115.8079
$ws.Cells.Item(132, 13).Value = -8585.8079
# Row 134: Guaranteed Gem
$ws.Cells.Item(134, 8).Value = 78053.37
$ws.Cells.Item(134, 10).Value = 78053.37
$ws.Cells.Item(134, 12).Value = 234160.11
$ws.Cells.Item(134, 14).Value = -239230.11

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Cells.Item(46, 8).Value = 6582.9585
$ws.Cells.Item(46, 10).Value = 10356.929
$ws.Cells.Item(46, 12).Value = 10356.929
$ws.Cells.Item(46, 14).Value = -10732.929
# Row 55: It's Not a Job, It's a Calling
$ws.Cells.Item(55, 8).Value = 642.7273
$ws.Cells.Item(55, 10).Value = 474.25
$ws.Cells.Item(55, 12).Value = 474.25
$ws.Cells.Item(55, 14).Value = -820.25
# Row 63: From Mud to Mourning
$ws.Cells.Item(63, 8).Value = 44999
$ws.Cells.Item(63, 9).Value = 44999
$ws.Cells.Item(63, 11).Value = 44999
$ws.Cells.Item(63, 13).Value = -44250
# Row 66: These Boots Are Made for Hawkin' (L)
$ws.Cells.Item(66, 8).Value = 44999
$ws.Cells.Item(66, 9).Value = 44999
$ws.Cells.Item(66, 11).Value = 134997
$ws.Cells.Item(66, 13).Value = -131253
# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 5819.2656
$ws.Cells.Item(132, 9).Value = 4507.5557
$ws.Cells.Item(132, 11).Value = 13522.6671
$ws.Cells.Item(132, 13).Value = -10992.6671
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 5606.4136
$ws.Cells.Item(136, 9).Value = 4322.8086
$ws.Cells.Item(136, 11).Value = 12968.4258
$ws.Cells.Item(136, 13).Value = -10418.4258

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Cells.Item(122, 8).Value = 5282.067
$ws.Cells.Item(122, 9).Value = 5282.067
$ws.Cells.Item(122, 11).Value = 15846.201
$ws.Cells.Item(122, 13).Value = -13396.201
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 3732.9768
$ws.Cells.Item(132, 10).Value = 5036.077
$ws.Cells.Item(132, 12).Value = 15108.231
$ws.Cells.Item(132, 14).Value = -20168.231
# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 1930.4242
$ws.Cells.Item(136, 9).Value = 1196.6471
$ws.Cells.Item(136, 11).Value = 3589.9413
$ws.Cells.Item(136, 13).Value = -1039.9413

